# Applies the cryptos list price/volume refresh described in the commit
# 'Updated cryptos list on Wed Aug 23 06:57:49 UTC 2023 with GitHub Actions'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look numeric (e.g. '1.001', '0.01620').
# Force the whole data range to Text format first so Excel does not coerce
# these into floating point numbers and strip significant trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.128.57'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '1.647.93'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '215.38'
$ws.Range('E5').Value = '  +2.77%  '
$ws.Range('D6').Value = '0.5233'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').Value = '0.2611'
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('D9').Value = '0.06325'
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('D10').Value = '20.83'
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('D11').Value = '0.07693'
$ws.Range('E11').Value = '  +2.15%  '
$ws.Range('D12').Value = '1.643.48'
$ws.Range('E12').Value = '  -1.25%  '
$ws.Range('E13').Value = '  -0.13%  '
$ws.Range('D14').Value = '1.862.22'
$ws.Range('E14').Value = '  -1.75%  '
$ws.Range('E15').Value = '  +1.73%  '
$ws.Range('D16').Value = '0.0₅8208'
$ws.Range('E16').Value = '  +3.21%  '
$ws.Range('D17').Value = '65.25'
$ws.Range('E17').Value = '  -1.56%  '
$ws.Range('D18').Value = '26.151.40'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('E20').Value = '  +0.36%  '
$ws.Range('D21').Value = '189.31'
$ws.Range('E21').Value = '  +1.56%  '
$ws.Range('D22').Value = '10.25'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').Value = '6.208'
$ws.Range('E23').Value = '  +0.52%  '
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('D25').Value = '145.69'
$ws.Range('E25').Value = '  -2.37%  '
$ws.Range('D26').Value = '7.440'
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('D27').Value = '0.1210'
$ws.Range('E27').Value = '  -3.00%  '
$ws.Range('D28').Value = '15.93'
$ws.Range('E28').Value = '  +0.71%  '
$ws.Range('D29').Value = '1.395'
$ws.Range('E29').Value = '  +3.18%  '
$ws.Range('D30').Value = '0.05891'
$ws.Range('E30').Value = '  -7.69%  '
$ws.Range('D31').Value = '1.262'
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').Value = '3.447'
$ws.Range('E32').Value = '  -1.40%  '
$ws.Range('D33').Value = '3.413'
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('D34').Value = '1.657'
$ws.Range('E34').Value = '  +0.87%  '
$ws.Range('D35').Value = '0.9868'
$ws.Range('E35').Value = '  -1.65%  '
$ws.Range('D36').Value = '2.768'
$ws.Range('E36').Value = '  +0.87%  '
$ws.Range('D37').Value = '2.390'
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('D38').Value = '0.5685'
$ws.Range('E38').Value = '  -5.30%  '
$ws.Range('D39').Value = '0.01620'
$ws.Range('E39').Value = '  +0.41%  '
$ws.Range('D40').Value = '0.8581'
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '5.759'
$ws.Range('E41').Value = '  -6.05%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '1.001'
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('D43').Value = '1.031.18'
$ws.Range('E43').Value = '  -7.14%  '
$ws.Range('D44').Value = '100.23'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').Value = '1.794.04'
$ws.Range('E45').Value = '  -1.36%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '56.04'
$ws.Range('E46').Value = '  +1.43%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₈105'
$ws.Range('E47').Value = '  -4.75%  '
$ws.Range('D48').Value = '1.004'
$ws.Range('E48').Value = '  +0.72%  '
$ws.Range('D49').Value = '8.125'
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('D50').Value = '0.05160'
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('D51').Value = '0.4222'
$ws.Range('E51').Value = '  -0.53%  '
